$d = $word.ActiveDocument
$p1 = $d.Paragraphs.Item(1)
$full = $p1.Range
$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w14:paraId="37E4074D" w14:textId="1E88BD9B" w:rsidR="001C6A97" w:rsidRPr="00B50BFD" w:rsidRDefault="00B50BFD" w:rsidP="00B50BFD"><w:pPr><w:spacing w:line="276" w:lineRule="auto"/><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Our class diagram begins with the Login class, which handles w</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">here the user will be directed </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>to</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>:</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">the manager, customer, wait staff, </w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">or </w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>kitchen staff</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> views.</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> From there, the manager and customer have access to the menu, allowing the manager to make changes to it and customers to view and add items to the cart. The manager will also have exclusive access to store information. From the Menu class, Order will have a dependency on Order, as Order cannot exist without Menu. The wait</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>staff, kitchen staff, and customer will all have access to the Order class. After the order the customer can play Games, so Games branches off Order. Lastly, Payment will come after order, so it is attached to that. The manager will also have access to the Payment class, as there are special cases where he might need to intervene with the payment procedure (compensating for an item).</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$full.InsertXML($xml)

# InsertXML on the full paragraph range (including its paragraph mark) leaves behind an
# extra blank paragraph after the newly inserted one; merge it back out.
if ($d.Paragraphs.Count -gt 1) {
  $np1 = $d.Paragraphs.Item(1)
  $np2 = $d.Paragraphs.Item(2)
  $mergeRange = $d.Range($np1.Range.End - 1, $np2.Range.End)
  $mergeRange.Delete()
}

Write-Output $d.Paragraphs.Item(1).Range.Text
